# Fixed bug for Spencer's Method for right-facing slopes
#
# The columns BC:BH (Fh, Fv, Mo, yt_l, yt_r -- plus the stray leftover
# `y_q`) are removed; BA/BB are renamed to yt_l/yt_r and recomputed, and
# the dependent n_eff/z/theta columns (AH/AI/AJ) are recalculated with the
# corrected (right-facing) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the now-unused trailing columns BC:BH (Q, y_q, Fh, Fv, Mo and
#    the old yt_l/yt_r position) so the used range shrinks to A1:BB22.
# ---------------------------------------------------------------------
$ws.Range("BC1:BH22").Clear()

# ---------------------------------------------------------------------
# 2. Re-header the two surviving columns.
# ---------------------------------------------------------------------
$ws.Range("BA1").Value = "yt_l"
$ws.Range("BB1").Value = "yt_r"

# ---------------------------------------------------------------------
# 3. Recomputed data, rows 2..22 (slice 1..21).
# ---------------------------------------------------------------------
$AH = @(-263.490129717897,4332.154894227278,29608.58027461288,27922.04264876728,33698.57345613423,38871.0276108088,43649.98656765717,48042.70950462885,52500.85350505389,47579.46138868868,44131.66947491431,40423.3848287001,44278.52931612825,39308.26420413349,35636.05959501668,28870.8861214081,29311.05773992906,29300.73188562266,28498.47933690739,26047.55986138142,19153.65912806178)
$AI = @(0,-351.0935560601921,10479.5238595941,71500.27834923461,114388.3728871403,154547.120363853,189467.6940758316,217653.0421241137,238115.2714811733,250022.4030152136,253026.7317882518,249443.7370002513,240170.5697123781,222622.4430803917,199782.3330560254,173310.8829600393,148008.1418949654,119112.5078362852,87441.31792979216,54461.40285713836,23042.20337559174)
$AJ = @(-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689,-8.895705300300689)
$BA = @(84.00000000000006,81.13625166065518,69.67467870263968,51.7531074920696,44.05438597026956,37.87874927014892,32.8925908389824,28.87466999988604,25.69169069832504,23.28078639983498,21.74951749706717,20.75245860086207,20.23640734131954,20.217882391099,20.77587039197478,21.82653761969504,23.13302111351473,25.02146089394151,27.58051703355543,30.95233001688217,35.40375914414403)
$BB = @(81.13625166065518,69.67467870263968,51.7531074920696,44.05438597026956,37.87874927014892,32.8925908389824,28.87466999988604,25.69169069832504,23.28078639983498,21.74951749706717,20.75245860086207,20.23640734131954,20.217882391099,20.77587039197478,21.82653761969504,23.13302111351473,25.02146089394151,27.58051703355543,30.95233001688217,35.40375914414403,39.9962422235724)

for ($i = 0; $i -lt 21; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 34).Value = $AH[$i]   # column AH
    $ws.Cells.Item($r, 35).Value = $AI[$i]   # column AI
    $ws.Cells.Item($r, 36).Value = $AJ[$i]   # column AJ
    $ws.Cells.Item($r, 53).Value = $BA[$i]   # column BA
    $ws.Cells.Item($r, 54).Value = $BB[$i]   # column BB
}
